$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "298.36"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.62%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.15%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.159"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.02%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.682"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "60.01%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.822"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.40%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.826"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.45%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9143"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.49%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1734"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.66%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07245"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.86%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08382"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.00%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03000"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.72%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09952"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.48%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001496"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.36%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006111"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.94%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.250"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.09%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.36%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1338"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.51%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.631"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.73%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04565"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.41%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001258"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.44%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004446"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.38%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-9.06%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003431"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "83.35%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01824"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "6.77%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04505"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.85%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007012"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1341"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.76%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002241"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.86%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009823"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.49%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006468"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.69%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.03%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.006202"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-39.27%"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-56.08%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.03%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.10%"
